$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new data row at row 164, pushing the existing rows
# 164..257 down to 165..258 (dimension grows from A1:R257 to A1:R258).
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(164, 1).Value = 10
$ws.Cells.Item(164, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(164, 3).Value = "La Araucanía"
$ws.Cells.Item(164, 4).Value = 44767
$ws.Cells.Item(164, 5).Value = 9
$ws.Cells.Item(164, 6).Value = 100112043
$ws.Cells.Item(164, 7).Value = "Pepino dulce"
$ws.Cells.Item(164, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 240
$ws.Cells.Item(164, 11).Value = 18000
$ws.Cells.Item(164, 12).Value = 19000
$ws.Cells.Item(164, 13).Value = 18583
$ws.Cells.Item(164, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(164, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(164, 16).Value = 1032
$ws.Cells.Item(164, 17).Value = 18
$ws.Cells.Item(164, 18).Value = "Hortaliza"
